# Fill in the missing "scaledTanh" / "Last Layer Only" results
# (row 6, columns G:I) on Sheet1 of the activation-function test results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("G6").Value = 4.5423
$ws.Range("H6").Value = 6.9443799999999998
$ws.Range("I6").Value = 6.2509300000000003

# Leave the freshly entered cell selected, matching the saved view state.
$ws.Range("I6").Select()
